$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Total Number of Lipids (updated counts)
$ws.Range("B2").Value = 579
$ws.Range("C2").Value = 388
$ws.Range("D2").Value = 799
$ws.Range("E2").Value = 540
$ws.Range("F2").Value = 603
$ws.Range("G2").Value = 353

# Row 3 - "<= 10" percentages
$ws.Range("B3").Value = "145 (25%)"
$ws.Range("C3").Value = "89 (23%)"
$ws.Range("D3").Value = "263 (33%)"
$ws.Range("F3").Value = "121 (20%)"
$ws.Range("G3").Value = "52(15%)"

# Row 4 - "<= 20" percentages
$ws.Range("B4").Value = "311 (54%)"
$ws.Range("C4").Value = "208 (54%)"
$ws.Range("D4").Value = "473 (59%)"
$ws.Range("F4").Value = "263 (44%)"
$ws.Range("G4").Value = "162 (46%)"

# Row 5 - "<= 30" percentages
$ws.Range("B5").Value = "424 (73%)"
$ws.Range("C5").Value = "293 (76%)"
$ws.Range("D5").Value = "599 (75%)"
$ws.Range("F5").Value = "418 (69%)"
$ws.Range("G5").Value = "259 (73%)"

# Row 6 - Remained Lipids counts (row 7 totals are formulas and recalc automatically)
$ws.Range("B6").Value = 424
$ws.Range("C6").Value = 293
$ws.Range("D6").Value = 599
$ws.Range("F6").Value = 418
$ws.Range("G6").Value = 259

# Update the active selection to match the final saved state
$ws.Range("F26").Select()
